$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(94, '2026-02-13 12:28:50', '237671825253', 'MAFFO YEMDJI CHIMENE ETS MOBILE FINANCIAL SERVICES MFS', 27439),
    @(95, '2026-02-13 17:13:39', '237673816350', 'LA NEGRESSE SARL DJUIDJE EPOUSE TAGNE HELENE ADELE', 102856),
    @(96, '2026-02-13 16:41:16', '237674484736', 'MAMADOU ALPHA DIALLO', 175549),
    @(97, '2026-02-13 15:52:43', '237674895525', 'NFOR EPSE FOMUNGUM ASSUMPTA MUMBEB', 287848),
    @(98, '2026-02-13 16:59:01', '237675457527', 'TCHOUANKAP DJOMKAM ISMAEL ETS MOBILE FINANCIAL SERVICES MFS', 121749),
    @(99, '2026-02-13 15:57:51', '237675551814', 'PELAGIE AIMEE NTOUBA MPAKO', 299530),
    @(100, '2026-02-13 12:33:12', '237678530662', 'Seraphine Abela Eyele', 122082),
    @(101, '2026-02-13 15:26:26', '237678796497', 'MANIGANG NDALLOKA MARIE JOSIANE ALBARKA GN SARL', 136760),
    @(102, '2026-02-13 18:23:22', '237679093371', 'ADAMA SARRE', 29745),
    @(103, '2026-02-09 08:26:36', '237679209479', 'TEDJON CLAUVIS FRANCIS TOP MOBIL TELECOM', 571),
    @(104, '2026-02-13 16:12:46', '237679732169', 'HILAIRE EBWANGA FOTSO', 76217),
    @(105, '2026-02-13 15:19:47', '237679791269', 'ETS LE CONTENT 67', 0),
    @(106, '2026-02-13 14:43:52', '237681663743', 'LA NEGRESSE SARL FONGA SINTCHA YOLANDE MIREILLE', 514278),
    @(107, '2026-02-13 13:31:20', '237682480811', 'KENGNE TADJO LYNDA NOEL ETS MOBILE FINANCIAL SERVICES MFS', 764658),
    @(108, '2026-02-13 11:57:18', '237682798275', 'NGAFFO YOCADINE BENEDITE ETS MOBILE FINANCIAL SERVICES MFS', 200036),
    @(109, '2026-02-13 16:11:27', '237683165199', 'N A ABASS GONI', 9058),
    @(110, '2026-02-13 15:50:13', '237683454059', 'RTS BELLE HOLLANDAISE', 223585)
)

foreach ($row in $data) {
    $r = $row[0]
    $dateStr = $row[1]
    $numStr = $row[2]
    $name = $row[3]
    $balance = $row[4]

    $ws.Range("A$r").Value = $dateStr

    $ws.Range("B$r").NumberFormat = "@"
    $ws.Range("B$r").Value = $numStr
    $ws.Range("B$r").Style = "Normal"

    $ws.Range("C$r").Value = $name
    $ws.Range("D$r").Value = $balance
}
